$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 7422965940
$ws.Range("B4").Value = "Vodafone"

$ws.Range("A4").Select()
